# Applies the "Common: A lot of intereting stuff" change:
#  - Translations - Common: insert one new row (common.job.name.migrate / Migrace)
#  - Translations - Lab: insert one new row (lab.mixture.liquid.create.button / Namíchat)
#    plus append 10 new rows for the new "mixture preview / liquid create" translation keys
#  - Active sheet moves from "Translations - Common" to "Translations - Lab"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Translations - Common: insert a new translation row at row 38
# ---------------------------------------------------------------------------
$wsCommon = $wb.Worksheets.Item("Translations - Common")

$wsCommon.Rows.Item(38).Insert()
$wsCommon.Range("A38").Value = "cs"
$wsCommon.Range("B38").Value = "common.job.name.migrate"
$wsCommon.Range("C38").Value = "Migrace"

# ---------------------------------------------------------------------------
# Translations - Lab: insert a new translation row at row 123
# ---------------------------------------------------------------------------
$wsLab = $wb.Worksheets.Item("Translations - Lab")

$wsLab.Rows.Item(123).Insert()
$wsLab.Range("A123").Value = "cs"
$wsLab.Range("B123").Value = "lab.mixture.liquid.create.button"
$wsLab.Range("C123").Value = "Namíchat"

# Append new rows with the newly added "mixture preview" / "liquid create" translations
# (copy the formatting of the last existing data row down onto the new rows first)
$wsLab.Range("A131:C131").Copy()
$wsLab.Range("A132:C141").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsLab.Range("A132").Value = "cs"
$wsLab.Range("B132").Value = "lab.mixture.liquid.create.title"
$wsLab.Range("C132").Value = "Nový liquid"

$wsLab.Range("A133").Value = "cs"
$wsLab.Range("B133").Value = "lab.mixture.liquid.create.content"
$wsLab.Range("C133").Value = "Přejete si namíchat liquid z vybraného mixu? Datum zrání započne jeho vytvořením."

$wsLab.Range("A134").Value = "cs"
$wsLab.Range("B134").Value = "lab.mixture.preview.aroma"
$wsLab.Range("C134").Value = "Aroma"

$wsLab.Range("A135").Value = "cs"
$wsLab.Range("B135").Value = "lab.mixture.preview.pgvg"
$wsLab.Range("C135").Value = "Poměr VG/PG"

$wsLab.Range("A136").Value = "cs"
$wsLab.Range("B136").Value = "lab.mixture.preview.content"
$wsLab.Range("C136").Value = "Obsah aromatu"

$wsLab.Range("A137").Value = "cs"
$wsLab.Range("B137").Value = "lab.mixture.preview.base"
$wsLab.Range("C137").Value = "Báze"

$wsLab.Range("A138").Value = "cs"
$wsLab.Range("B138").Value = "lab.mixture.preview.booster"
$wsLab.Range("C138").Value = "Booster"

$wsLab.Range("A139").Value = "cs"
$wsLab.Range("B139").Value = "lab.mixture.liquid.create.ok.button"
$wsLab.Range("C139").Value = "Namíchat"

$wsLab.Range("A140").Value = "cs"
$wsLab.Range("B140").Value = "lab.mixture.liquid.create.success"
$wsLab.Range("C140").Value = "Liquid úspěšně namíchán."

$wsLab.Range("A141").Value = "cs"
$wsLab.Range("B141").Value = "lab.mixture.liquid.create.failure"
$wsLab.Range("C141").Value = "Liquid se nepodařilo namíchat."

# ---------------------------------------------------------------------------
# Selections on both sheets (matches the recorded view state of the edit)
# ---------------------------------------------------------------------------
$wsCommon.Range("B65").Select()

# Activate "Translations - Lab" last so it becomes the workbook's active tab
$wsLab.Activate()
$wsLab.Range("B135").Select()
